$wb = $excel.ActiveWorkbook

# =========================================================================
# "Generate Report for Handoff" - the CI report is regenerated: rows are
# reordered (the handed-off-but-not-yet-handed-back file moves to the
# bottom with a refreshed "Ready for handoff" status/time) and some
# timestamps / generated-file references are refreshed.
# Only the cells whose resolved text actually changes are touched below.
# =========================================================================

# ---- Overview ----
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md"
$ws.Range("B2").Value = "e2e\ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md"
$ws.Range("G2").Value = "2016-08-19 02:58:56"
$ws.Range("A3").Value = "ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md"
$ws.Range("B3").Value = "e2e\ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md"
$ws.Range("A4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md"
$ws.Range("B4").Value = "e2e\7a3bce59-30f1-4a28-96a7-2f3be1a63706.md"
$ws.Range("E4").Value = "Ready for handoff"
$ws.Range("F4").Value = "Ready for handoff"
$ws.Range("G4").Value = "2016-08-19 03:01:58"

# rebuild hyperlinks (display text per column must follow the new row
# order while the underlying link targets stay attached to the same
# B2/B3/B4 positions)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b9fd50b174caf4b40afb9c1b61e649f4ca8b549/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md", [Type]::Missing, [Type]::Missing, "e2e\ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/102353cc41196c4c326f61553dd2802c3da2660d/e2e/ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md", [Type]::Missing, [Type]::Missing, "e2e\ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b9fd50b174caf4b40afb9c1b61e649f4ca8b549/e2e/ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md", [Type]::Missing, [Type]::Missing, "e2e\7a3bce59-30f1-4a28-96a7-2f3be1a63706.md") | Out-Null

# ---- zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md"
$ws.Range("G2").Value = "d64ea6bd-8584-404b-89fa-11aa8c748f12.7d3960ec2ad0b4f3138cea710972e3a654c9f604.zh-cn.xlf"
$ws.Range("H2").Value = "2016-08-19 02:58:51"
$ws.Range("I2").Value = "d64ea6bd-8584-404b-89fa-11aa8c748f12.md"
$ws.Range("J2").Value = "d64ea6bd-8584-404b-89fa-11aa8c748f12.7d3960ec2ad0b4f3138cea710972e3a654c9f604.zh-cn.xlf"
$ws.Range("K2").Value = "2016-08-19 02:59:15"
$ws.Range("A3").Value = "ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md"
$ws.Range("F3").Value = "True"
$ws.Range("A4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "False"
$ws.Range("G4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.b86d21e752e43d2bcf73a5fe10355cacf5943330.zh-cn.xlf"
$ws.Range("H4").Value = "2016-08-19 03:01:52"
$ws.Range("I4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md"
$ws.Range("J4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.b86d21e752e43d2bcf73a5fe10355cacf5943330.zh-cn.xlf"
$ws.Range("K4").Value = "2016-08-19 03:01:28"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b9fd50b174caf4b40afb9c1b61e649f4ca8b549/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3fa5806600ab640edba30c9089dfcacf530f385/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md."

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b9fd50b174caf4b40afb9c1b61e649f4ca8b549/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md", [Type]::Missing, [Type]::Missing, "ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e4a75490e3f340641f00a244bab9be531d6b1ccd/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md", [Type]::Missing, [Type]::Missing, "d64ea6bd-8584-404b-89fa-11aa8c748f12.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/102353cc41196c4c326f61553dd2802c3da2660d/e2e/ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md", [Type]::Missing, [Type]::Missing, "ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/91e601ae424690594fbf56d88b15052da6e43a0e/e2e/d64ea6bd-8584-404b-89fa-11aa8c748f12.md", [Type]::Missing, [Type]::Missing, "d64ea6bd-8584-404b-89fa-11aa8c748f12.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b9fd50b174caf4b40afb9c1b61e649f4ca8b549/e2e/ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md", [Type]::Missing, [Type]::Missing, "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/91e601ae424690594fbf56d88b15052da6e43a0e/e2e/d64ea6bd-8584-404b-89fa-11aa8c748f12.md", [Type]::Missing, [Type]::Missing, "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md") | Out-Null

# column P now holds the long "not latest" message -> widen it like the source report does
$ws.Columns.Item(16).ColumnWidth = 39.15

# ---- de-de ----
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md"
$ws.Range("G2").Value = "d64ea6bd-8584-404b-89fa-11aa8c748f12.7d3960ec2ad0b4f3138cea710972e3a654c9f604.de-de.xlf"
$ws.Range("H2").Value = "2016-08-19 02:58:56"
$ws.Range("I2").Value = "d64ea6bd-8584-404b-89fa-11aa8c748f12.md"
$ws.Range("J2").Value = "d64ea6bd-8584-404b-89fa-11aa8c748f12.7d3960ec2ad0b4f3138cea710972e3a654c9f604.de-de.xlf"
$ws.Range("K2").Value = "2016-08-19 02:59:22"
$ws.Range("A3").Value = "ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md"
$ws.Range("F3").Value = "True"
$ws.Range("A4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("F4").Value = "False"
$ws.Range("G4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.b86d21e752e43d2bcf73a5fe10355cacf5943330.de-de.xlf"
$ws.Range("H4").Value = "2016-08-19 03:01:58"
$ws.Range("I4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md"
$ws.Range("J4").Value = "7a3bce59-30f1-4a28-96a7-2f3be1a63706.b86d21e752e43d2bcf73a5fe10355cacf5943330.de-de.xlf"
$ws.Range("K4").Value = "2016-08-19 03:01:35"
$ws.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b9fd50b174caf4b40afb9c1b61e649f4ca8b549/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3fa5806600ab640edba30c9089dfcacf530f385/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md."

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b9fd50b174caf4b40afb9c1b61e649f4ca8b549/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md", [Type]::Missing, [Type]::Missing, "ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/101e58650b15b54caa990151a6e303476c1cd08c/e2e/7a3bce59-30f1-4a28-96a7-2f3be1a63706.md", [Type]::Missing, [Type]::Missing, "d64ea6bd-8584-404b-89fa-11aa8c748f12.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/102353cc41196c4c326f61553dd2802c3da2660d/e2e/ffff1065edc8-172e-4b1c-9000-0bf1e391a5df.md", [Type]::Missing, [Type]::Missing, "ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2456f20ad0ec5dbeabddd98a8d58bc570c3dd443/e2e/d64ea6bd-8584-404b-89fa-11aa8c748f12.md", [Type]::Missing, [Type]::Missing, "d64ea6bd-8584-404b-89fa-11aa8c748f12.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2b9fd50b174caf4b40afb9c1b61e649f4ca8b549/e2e/ffffff18d0d5ab-13b4-4127-9ec4-418b0643280b.md", [Type]::Missing, [Type]::Missing, "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2456f20ad0ec5dbeabddd98a8d58bc570c3dd443/e2e/d64ea6bd-8584-404b-89fa-11aa8c748f12.md", [Type]::Missing, [Type]::Missing, "7a3bce59-30f1-4a28-96a7-2f3be1a63706.md") | Out-Null

# column P now holds the long "not latest" message -> widen it like the source report does
$ws.Columns.Item(16).ColumnWidth = 39.15

$wb.Save()
